$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the 8 new rows (5-12) by copying row 4 down, so border/alignment
#     styling on column A (and row structure) matches the existing rows. ---
$ws.Range("A4:H4").Copy() | Out-Null
$ws.Range("A5:H12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Sequential index values in column A (0-based task order)
for ($r = 2; $r -le 12; $r++) {
    $ws.Range("A$r").Value = ($r - 2)
}

# --- Column D/E for rows 4-12 hold numeric-looking values that must be
#     stored as TEXT (shared strings) rather than numbers; rows 2-3 stay numeric. ---
$ws.Range("D4:E12").NumberFormat = "@"

# Row 2
$ws.Range("B2").Value = "Virtual Desk Timeout"
$ws.Range("C2").Value = "NonCode"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = "2018-08-09`n10:00`n2 days left!"
$ws.Range("G2").Value = "work/misc"
$ws.Range("H2").Value = "None"

# Row 3
$ws.Range("B3").Value = "Data Transfer Rate"
$ws.Range("C3").Value = "NoneCode"
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = "2018-08-10`n12:00`n3 days left!"
$ws.Range("G3").Value = "work/misc"
$ws.Range("H3").Value = "None"

# Row 4
$ws.Range("B4").Value = "Better QC Tool"
$ws.Range("C4").Value = "Code"
$ws.Range("D4").Value = "6"
$ws.Range("E4").Value = "7"
$ws.Range("F4").Value = "2018-08-11`n14:00`n5 days and 0 hours left!"
$ws.Range("G4").Value = "work/misc"
$ws.Range("H4").Value = "None"

# Row 5
$ws.Range("B5").Value = "Automate Signout Checklist"
$ws.Range("C5").Value = "Code"
$ws.Range("D5").Value = "6"
$ws.Range("E5").Value = "8"
$ws.Range("F5").Value = "2018-08-12`n08:00`n5 days and 18 hours left!"
$ws.Range("G5").Value = "work/misc"
$ws.Range("H5").Value = "None"

# Row 6
$ws.Range("B6").Value = "Increase CPU Power"
$ws.Range("C6").Value = "NonCode"
$ws.Range("D6").Value = "3"
$ws.Range("E6").Value = "3"
$ws.Range("F6").Value = "2018-08-13`n23:00`n7 days and 9 hours left!"
$ws.Range("G6").Value = "work/misc"
$ws.Range("H6").Value = "None"

# Row 7
$ws.Range("B7").Value = "ngjob file sample order"
$ws.Range("C7").Value = "NonCode"
$ws.Range("D7").Value = "3"
$ws.Range("E7").Value = "9"
$ws.Range("F7").Value = "2018-08-15`n03:00`n8 days and 13 hours left!"
$ws.Range("G7").Value = "work/misc"
$ws.Range("H7").Value = "None"

# Row 8
$ws.Range("B8").Value = "Primer Seek"
$ws.Range("C8").Value = "Code"
$ws.Range("D8").Value = "10"
$ws.Range("E8").Value = "15"
$ws.Range("F8").Value = "2018-08-20`n12:00`n13 days and 22 hours left!"
$ws.Range("G8").Value = "work/misc"
$ws.Range("H8").Value = "None"

# Row 9
$ws.Range("B9").Value = "Automate Execution of NextGene"
$ws.Range("C9").Value = "Code"
$ws.Range("D9").Value = "5"
$ws.Range("E9").Value = "10"
$ws.Range("F9").Value = "2018-08-01`n10:00`nTask is past due"
$ws.Range("G9").Value = "work/misc"
$ws.Range("H9").Value = "None"

# Row 10
$ws.Range("B10").Value = "BAMasker/BRR One-click"
$ws.Range("C10").Value = "Code"
$ws.Range("D10").Value = "6"
$ws.Range("E10").Value = "6"
$ws.Range("F10").Value = "2018-08-05`n00:00`nTask is past due"
$ws.Range("G10").Value = "work/misc"
$ws.Range("H10").Value = "None"

# Row 11
$ws.Range("B11").Value = "Automate data transfer"
$ws.Range("C11").Value = "Code"
$ws.Range("D11").Value = "7"
$ws.Range("E11").Value = "12"
$ws.Range("F11").Value = "2018-08-06`n10:00`nTask is past due"
$ws.Range("G11").Value = "work/misc"
$ws.Range("H11").Value = "None"

# Row 12
$ws.Range("B12").Value = "Automate VCF Upload Annotation"
$ws.Range("C12").Value = "Code"
$ws.Range("D12").Value = "8"
$ws.Range("E12").Value = "7"
$ws.Range("F12").Value = "2018-08-06`n17:00`nDue today!"
$ws.Range("G12").Value = "work/misc"
$ws.Range("H12").Value = "None"

# Multi-line Deadline (column F) values trigger an automatic custom row height;
# auto-fit each affected row back down so no explicit row height is stored.
# (Must run before ClearFormats below so the text number-format does not get
#  re-stamped onto D:E after the height recalculation.)
$ws.Rows("2:12").AutoFit() | Out-Null

# Cells written above as text picked up an implicit "@" text style; clear it so
# they fall back to the default (unstyled) cell format, matching the target file.
$ws.Range("D4:E12").ClearFormats()
